$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 1085

$data = @(
    @(45534.5,538.4,543.4,522.2,523.9,77675.621),
    @(45534.66666666666,523.8,537.1,522,531.6,35832.741),
    @(45534.83333333334,531.6,536.4,531,535.2,9566.069),
    @(45535,535.2,538,534.8,536,6077.518),
    @(45535.16666666666,537.4,541.5,535.1,535.6,24767.867),
    @(45535.33333333334,535.7,537.7,534.2,536.6,13144.509),
    @(45535.5,536.6,537.8,535,535.8,17108.452),
    @(45535.66666666666,535.8,535.9,531.1,533.5,13539.582),
    @(45535.83333333334,533.6,533.9,531.6,532.7,4273.418),
    @(45536,532.9,533.2,524.7,526.8,27867.97),
    @(45536.16666666666,526.8,527,513.4,522.3,63164.734),
    @(45536.33333333334,522.3,522.6,516.1,519,17177.512),
    @(45536.5,519.1,522,510.3,520.9,42419.228),
    @(45536.66666666666,520.9,524,516.1,520.2,23621.356),
    @(45536.83333333334,520.3,523.3,507.7,512.4,29448.067),
    @(45537,512.3,516.5,511.4,515.2,18951.789),
    @(45537.16666666666,515.1,515.9,501.5,503.1,80401.34),
    @(45537.33333333334,503.1,525.4,502.9,519.3,98970.88),
    @(45537.5,519.1,520.1,517.6,518.3,11262.374),
    @(45537.66666666666,522.3,523.9,519.8,522.2,13011.087),
    @(45537.83333333334,522.3,524.2,521.9,523.9,1007.043),
    @(45538,526.4,538.3,525.2,534.1,52774.763),
    @(45538.16666666666,534.1,538.4,531.6,535,46945.48),
    @(45538.33333333334,535,536.2,531.6,535.9,22313.805),
    @(45538.5,535.9,537.9,521.2,522,52347.766),
    @(45538.66666666666,522,525.7,519.6,524.4,25841.588),
    @(45538.83333333334,524.4,526.1,517.9,518.7,24258.248),
    @(45539,518.8,523.7,501.3,512.6,129095.058),
    @(45539.16666666666,512.7,513.7,508.3,510.2,34634.842),
    @(45539.33333333334,510.2,512.3,496.9,499,77841.611),
    @(45539.5,498.9,513.4,495.8,508.9,125237.795),
    @(45539.66666666666,509,512.6,505.3,509.6,45237.095),
    @(45539.83333333334,509.6,513.5,507,507.6,23072.62),
    @(45540,507.6,512.9,503.4,504,39281.877),
    @(45540.16666666666,504.1,508,501,506.7,30267.673),
    @(45540.33333333334,506.7,507.7,502.6,504.5,15940.201),
    @(45540.5,504.5,510.3,499.7,501.3,47496.706),
    @(45540.66666666666,501.4,503.6,496.7,499.5,32040.449),
    @(45540.83333333334,499.4,501.6,498.1,499.9,8542.011),
    @(45541,502.4,508,500.6,507.2,17860.734),
    @(45541.16666666666,505.4,508.4,496.9,501.4,37195.682),
    @(45541.33333333334,501.4,506.8,500.7,504.1,31446.474),
    @(45541.5,504.2,504.7,504.1,504.6,43.611)
)

$i = 0
foreach ($row in $data) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value2 = $row[0]
    $ws.Cells.Item($r, 2).Value2 = $row[1]
    $ws.Cells.Item($r, 3).Value2 = $row[2]
    $ws.Cells.Item($r, 4).Value2 = $row[3]
    $ws.Cells.Item($r, 5).Value2 = $row[4]
    $ws.Cells.Item($r, 6).Value2 = $row[5]
    $i++
}

# Copy style (border/alignment/font) of existing date column cell into new date cells
$srcStyleCell = $ws.Range("A1084")
$destRange = $ws.Range("A1085:A1127")
$srcStyleCell.Copy()
$destRange.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
